$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "showing all available WoRMs images" -- every species-name cell (A1:A26) gets
# wrap-text formatting turned on (direct formatting layered on top of the
# existing "XLConnect.String" cell style) so that the full name is visible once
# rows are resized to show the associated WoRMs thumbnail images.
$lastRow = 26
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.WrapText = $true
}
